# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on sheet "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.97 = 32664.53 pesos`n✅ 32664.53 pesos = 7.93 = 955.16 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the exchange-rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 125.399
$wsTasas.Range("O10").Value = 4096.1

$wsTasas.Range("N12").Value = 4119.18
$wsTasas.Range("O12").Value = 120.451
